# Judge validations - update sample judge row (rows 8-9) from
# "Anjali Abshire" / "BVAAABSHIRE" to "Roth, Lauren" / "DSUSER".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "Roth, Lauren"
$ws.Range("C8").Value = "DSUSER"
$ws.Range("B9").Value = "Roth, Lauren"
$ws.Range("C9").Value = "DSUSER"
